$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.302.31"
$ws.Range("E2").Value = "  +4.20%  "
$ws.Range("D3").Value = "1.809.76"
$ws.Range("E3").Value = "  +1.70%  "
$ws.Range("E4").Value = "  -0.68%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "338.78"
$ws.Range("E5").Value = "  +1.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9988"
$ws.Range("E6").Value = "  -0.61%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3925"
$ws.Range("E7").Value = "  +4.13%  "
$ws.Range("E8").Value = "  +2.66%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.59"
$ws.Range("E9").Value = "  +1.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.180"
$ws.Range("E10").Value = "  -0.50%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07539"
$ws.Range("E11").Value = "  +1.85%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9987"
$ws.Range("E12").Value = "  -0.72%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.05"
$ws.Range("E13").Value = "  +2.96%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.521"
$ws.Range("E14").Value = "  +2.41%  "
$ws.Range("D15").Value = "1.811.44"
$ws.Range("E15").Value = "  +1.74%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.161"
$ws.Range("E16").Value = "  +2.28%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001103"
$ws.Range("E17").Value = "  +2.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06706"
$ws.Range("E18").Value = "  +0.88%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "85.26"
$ws.Range("E19").Value = "  +1.82%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9999"
$ws.Range("E20").Value = "  -0.33%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.75"
$ws.Range("E21").Value = "  +3.66%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.561"
$ws.Range("E22").Value = "  +0.77%  "
$ws.Range("D23").Value = "28.300.58"
$ws.Range("E23").Value = "  +4.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.45"
$ws.Range("E24").Value = "  +0.91%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.400"
$ws.Range("E25").Value = "  -0.64%  "
$ws.Range("B26").Value = "ImmutableX"
$ws.Range("C26").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.476"
$ws.Range("E26").Value = "  -0.42%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.37"
$ws.Range("E27").Value = "  +1.98%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.521"
$ws.Range("E28").Value = "  +0.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "154.75"
$ws.Range("E29").Value = "  +1.85%  "
$ws.Range("D30").Value = "2.018.63"
$ws.Range("E30").Value = "  +1.86%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "136.12"
$ws.Range("E31").Value = "  +3.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.262"
$ws.Range("E32").Value = "  +5.28%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.019"
$ws.Range("E33").Value = "  -1.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08844"
$ws.Range("E34").Value = "  +2.60%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "13.16"
$ws.Range("E35").Value = "  +1.75%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02451"
$ws.Range("E36").Value = "  +5.53%  "
$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.470"
$ws.Range("E37").Value = "  +1.98%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06552"
$ws.Range("E38").Value = "  +4.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6923"
$ws.Range("E39").Value = "  +2.32%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.617"
$ws.Range("E40").Value = "  -1.82%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2230"
$ws.Range("E41").Value = "  +2.94%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.266"
$ws.Range("E42").Value = "  +1.88%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.563"
$ws.Range("E43").Value = "  -1.43%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.56"
$ws.Range("E44").Value = "  +1.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9985"
$ws.Range("E45").Value = "  -0.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6423"
$ws.Range("E46").Value = "  +1.89%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.874"
$ws.Range("E47").Value = "  +1.10%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.161"
$ws.Range("E48").Value = "  +2.70%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "131.74"
$ws.Range("E49").Value = "  +2.82%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07232"
$ws.Range("E50").Value = "  +1.08%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "80.34"
$ws.Range("E51").Value = "  +2.06%  "
